$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.906.68'
$ws.Range('E2').Value = '  +4.49%  '
$ws.Range('D3').Value = '3.391.05'
$ws.Range('E3').Value = '  +2.22%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '595.33'
$ws.Range('E5').Value = '  +6.95%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '186.99'
$ws.Range('E6').Value = '  +0.44%  '
$ws.Range('E7').Value = '  +4.21%  '
$ws.Range('E9').Value = '  +4.58%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.591'
$ws.Range('E10').Value = '  +2.23%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '47.68'
$ws.Range('E11').Value = '  +4.11%  '
$ws.Range('E12').Value = '  +6.80%  '
$ws.Range('D13').Value = '3.937.95'
$ws.Range('E13').Value = '  +2.70%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '639.95'
$ws.Range('E14').Value = '  +11.36%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '8.64'
$ws.Range('E15').Value = '  +2.19%  '
$ws.Range('D16').Value = '68.941.24'
$ws.Range('E16').Value = '  +4.63%  '
$ws.Range('D17').Value = '3.397.97'
$ws.Range('E17').Value = '  +3.05%  '
$ws.Range('E18').Value = '  +2.02%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '18.08'
$ws.Range('E19').Value = '  +2.11%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.17'
$ws.Range('E20').Value = '  +2.86%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.916'
$ws.Range('E21').Value = '  +2.77%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '18.03'
$ws.Range('E22').Value = '  -0.13%  '
$ws.Range('E23').Value = '  +2.45%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '100.26'
$ws.Range('E24').Value = '  +2.54%  '
$ws.Range('E25').Value = '  +3.74%  '
$ws.Range('E26').Value = '  +6.73%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.86'
$ws.Range('E27').Value = '  +5.15%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '32.95'
$ws.Range('E28').Value = '  +7.94%  '
$ws.Range('E29').Value = '  +3.69%  '
$ws.Range('E30').Value = '  +3.76%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '614.45'
$ws.Range('E31').Value = '  +9.23%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.79'
$ws.Range('E32').Value = '  +2.17%  '
$ws.Range('D33').Value = '4.028.64'
$ws.Range('E33').Value = '  +8.50%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '11.15'
$ws.Range('E34').Value = '  +3.06%  '
$ws.Range('E35').Value = '  +3.10%  '
$ws.Range('E36').Value = '  -0.19%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '56.70'
$ws.Range('E37').Value = '  +2.20%  '
$ws.Range('E38').Value = '  +8.23%  '
$ws.Range('B39').Value = 'Stacks'
$ws.Range('C39').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.34'
$ws.Range('E39').Value = '  +7.06%  '
$ws.Range('B40').Value = 'Kaspa'
$ws.Range('C40').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.131'
$ws.Range('E40').Value = '  +3.85%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '33.77'
$ws.Range('E41').Value = '  -0.45%  '
$ws.Range('D42').Value = '0.0₃0708'
$ws.Range('E42').Value = '  +3.21%  '
$ws.Range('E43').Value = '  +1.93%  '
$ws.Range('E44').Value = '  +3.50%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0425'
$ws.Range('E45').Value = '  +4.53%  '
$ws.Range('E46').Value = '  +1.95%  '
$ws.Range('E47').Value = '  +3.76%  '
$ws.Range('E48').Value = '  +12.67%  '
$ws.Range('E49').Value = '  +0.36%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '130.52'
$ws.Range('E50').Value = '  +2.68%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.86'
$ws.Range('E51').Value = '  +8.05%  '
